$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to Text format for the rows whose price
# value changes, so numeric-looking strings (e.g. "1.001") are stored
# as literal text instead of being coerced into numbers - matching the
# original inlineStr/text cells produced by the scraper.
$ws.Range("D2:D11").NumberFormat = "@"
$ws.Range("D13:D31").NumberFormat = "@"
$ws.Range("D33:D51").NumberFormat = "@"

# Per-row cell updates (Coin / Link / Price / Volume(1h))
# Row 2
$ws.Range("D2").Value = '29.152.88'
$ws.Range("E2").Value = '  -0.62%  '

# Row 3
$ws.Range("D3").Value = '1.823.60'
$ws.Range("E3").Value = '  -0.93%  '

# Row 4
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.23%  '

# Row 5
$ws.Range("D5").Value = '235.07'
$ws.Range("E5").Value = '  -1.68%  '

# Row 6
$ws.Range("D6").Value = '0.6063'
$ws.Range("E6").Value = '  -3.56%  '

# Row 7
$ws.Range("D7").Value = '1.001'
$ws.Range("E7").Value = '  +0.15%  '

# Row 8
$ws.Range("D8").Value = '0.07071'
$ws.Range("E8").Value = '  -4.85%  '

# Row 9
$ws.Range("D9").Value = '0.2795'
$ws.Range("E9").Value = '  -3.33%  '

# Row 10
$ws.Range("D10").Value = '23.45'
$ws.Range("E10").Value = '  -6.20%  '

# Row 11
$ws.Range("D11").Value = '0.07650'
$ws.Range("E11").Value = '  -0.95%  '

# Row 12
$ws.Range("E12").Value = '  -4.62%  '

# Row 13
$ws.Range("D13").Value = '4.789'
$ws.Range("E13").Value = '  -3.23%  '

# Row 14
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").Value = '0.6301'
$ws.Range("E14").Value = '  -6.49%  '

# Row 15
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").Value = '0.000009928'
$ws.Range("E15").Value = '  -2.44%  '

# Row 16
$ws.Range("D16").Value = '2.066.37'
$ws.Range("E16").Value = '  -1.00%  '

# Row 17
$ws.Range("D17").Value = '78.49'
$ws.Range("E17").Value = '  -3.83%  '

# Row 18
$ws.Range("D18").Value = '5.844'
$ws.Range("E18").Value = '  -5.96%  '

# Row 19
$ws.Range("D19").Value = '29.133.60'
$ws.Range("E19").Value = '  -0.51%  '

# Row 20
$ws.Range("D20").Value = '226.12'
$ws.Range("E20").Value = '  -1.27%  '

# Row 21
$ws.Range("D21").Value = '1.002'
$ws.Range("E21").Value = '  +0.19%  '

# Row 22
$ws.Range("D22").Value = '11.73'
$ws.Range("E22").Value = '  -4.45%  '

# Row 23
$ws.Range("D23").Value = '6.958'
$ws.Range("E23").Value = '  -5.15%  '

# Row 24
$ws.Range("D24").Value = '0.9999'
$ws.Range("E24").Value = '  +0.07%  '

# Row 25
$ws.Range("D25").Value = '155.15'
$ws.Range("E25").Value = '  -1.89%  '

# Row 26
$ws.Range("D26").Value = '8.014'
$ws.Range("E26").Value = '  -5.27%  '

# Row 27
$ws.Range("D27").Value = '0.1300'
$ws.Range("E27").Value = '  -3.29%  '

# Row 28
$ws.Range("D28").Value = '16.54'
$ws.Range("E28").Value = '  -4.64%  '

# Row 29
$ws.Range("D29").Value = '1.493'
$ws.Range("E29").Value = '  +2.47%  '

# Row 30
$ws.Range("D30").Value = '0.06247'
$ws.Range("E30").Value = '  -16.59%  '

# Row 31
$ws.Range("D31").Value = '1.445'
$ws.Range("E31").Value = '  -2.02%  '

# Row 32
$ws.Range("E32").Value = '  -5.23%  '

# Row 33
$ws.Range("D33").Value = '3.791'
$ws.Range("E33").Value = '  -6.08%  '

# Row 34
$ws.Range("D34").Value = '1.121'
$ws.Range("E34").Value = '  -1.54%  '

# Row 35
$ws.Range("D35").Value = '1.740'
$ws.Range("E35").Value = '  -4.11%  '

# Row 36
$ws.Range("D36").Value = '0.6378'
$ws.Range("E36").Value = '  -7.93%  '

# Row 37
$ws.Range("D37").Value = '2.543'
$ws.Range("E37").Value = '  -1.14%  '

# Row 38
$ws.Range("D38").Value = '1.211.03'
$ws.Range("E38").Value = '  -1.88%  '

# Row 39
$ws.Range("D39").Value = '2.714'
$ws.Range("E39").Value = '  -3.33%  '

# Row 40
$ws.Range("D40").Value = '0.01735'
$ws.Range("E40").Value = '  -5.59%  '

# Row 41
$ws.Range("D41").Value = '6.500'
$ws.Range("E41").Value = '  -5.39%  '

# Row 42
$ws.Range("D42").Value = '0.9034'
$ws.Range("E42").Value = '  -2.87%  '

# Row 43
$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  +0.11%  '

# Row 44
$ws.Range("B44").Value = 'RocketPoolETH'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D44").Value = '1.978.12'
$ws.Range("E44").Value = '  +0.37%  '

# Row 45
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").Value = '100.54'
$ws.Range("E45").Value = '  -0.32%  '

# Row 46
$ws.Range("D46").Value = '62.46'
$ws.Range("E46").Value = '  -4.26%  '

# Row 47
$ws.Range("D47").Value = '0.00000000115'
$ws.Range("E47").Value = '  -4.27%  '

# Row 48
$ws.Range("D48").Value = '1.593'
$ws.Range("E48").Value = '  -6.58%  '

# Row 49
$ws.Range("D49").Value = '8.471'
$ws.Range("E49").Value = '  -4.28%  '

# Row 50
$ws.Range("D50").Value = '0.4561'
$ws.Range("E50").Value = '  -0.53%  '

# Row 51
$ws.Range("D51").Value = '0.05505'
$ws.Range("E51").Value = '  -2.79%  '

